# DieselTrucks and total system Diesel Import MWh KPI
# - Rename the single-truck / scaled-vehicle asset in conversionAssets from
#   "TruckDiesel" / "VEHICLE" to "DieselTruck" / "DIESEL_VEHICLE".
# - Add the energy_consumption_kwhpkm / vehicle_scaling columns (J, K) to the
#   conversionAssets sheet, and populate them (plus the remaining numeric
#   columns) for the DieselTruck row.
# - Switch the active sheet/selection from storageAssets to conversionAssets.

$wb = $excel.ActiveWorkbook

$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsStorage    = $wb.Worksheets.Item("storageAssets")

# Rename the asset + type strings used by the vehicle-scaling diesel truck row.
$wsConversion.Range("B15").Value = "DieselTruck"
$wsConversion.Range("D15").Value = "DIESEL_VEHICLE"

# New header columns for energy consumption / vehicle scaling.
$wsConversion.Range("J1").Value = "energy_consumption_kwhpkm"
$wsConversion.Range("K1").Value = "vehicle_scaling"

# Fill in the rest of row 15 (DieselTruck) data.
$wsConversion.Range("E15").Value = 0
$wsConversion.Range("F15").Value = 0
$wsConversion.Range("G15").Value = 0.2
$wsConversion.Range("H15").Value = 0
$wsConversion.Range("J15").Value = 1

# storageAssets was the active/selected sheet before; conversionAssets becomes
# the active one now, with K15 selected. storageAssets keeps its own
# selection, now at O1.
$wsStorage.Activate()
$wsStorage.Range("O1").Select()

$wsConversion.Activate()
$wsConversion.Range("K15").Select()
